# Apply "Latest changes in framework" updates to the test-data workbook.
$wb = $excel.ActiveWorkbook

# --- DemoWebShop sheet: refresh generated Email / Firstname / Lastname for TC 1 ---
$wsDemoWebShop = $wb.Worksheets.Item("DemoWebShop")
$wsDemoWebShop.Range("C2").Value = "cYqcqcYb@gmail.com"
$wsDemoWebShop.Range("F2").Value = "BhoUpu"
$wsDemoWebShop.Range("G2").Value = "rSWvWf"

# --- OrangeHRM sheet: refresh generated Message / Middle Name / Last Name for TC 6 ---
$wsOrangeHRM = $wb.Worksheets.Item("OrangeHRM")
$wsOrangeHRM.Range("F6").Value = "thyRXv"
$wsOrangeHRM.Range("H6").Value = "thyRXv"
$wsOrangeHRM.Range("I6").Value = "thyRXv"
